{"js": "// Replace the 25 two-digit division prompts in the practice table with a\n// new set of problems, preserving each cell's existing run formatting\n// (font/size) and paragraph formatting (justification).\n//\n// The prompts live in table rows 0, 4, 8, 12 and 16 (5 cells each, one\n// problem per cell); the other rows are blank \"work space\" rows. Some of\n// the old values repeat (\"77\u00f72=\" appears twice) and some new values equal\n// an old value used elsewhere (\"36\u00f72=\", \"60\u00f75=\"), so replacements are\n// targeted at an exact (row, column) cell rather than found by a\n// document-wide text search - that keeps the edit unambiguous and order\n// independent.\nconst cellReplacements = [\n  // row 0\n  { row: 0, col: 0, oldText: \"79\u00f72=\", newText: \"22\u00f77=\" },\n  { row: 0, col: 1, oldText: \"56\u00f78=\", newText: \"62\u00f72=\" },\n  { row: 0, col: 2, oldText: \"55\u00f78=\", newText: \"28\u00f79=\" },\n  { row: 0, col: 3, oldText: \"43\u00f72=\", newText: \"73\u00f79=\" },\n  { row: 0, col: 4, oldText: \"84\u00f74=\", newText: \"40\u00f73=\" },\n  // row 4\n  { row: 4, col: 0, oldText: \"33\u00f78=\", newText: \"30\u00f74=\" },\n  { row: 4, col: 1, oldText: \"97\u00f73=\", newText: \"57\u00f72=\" },\n  { row: 4, col: 2, oldText: \"46\u00f74=\", newText: \"95\u00f75=\" },\n  { row: 4, col: 3, oldText: \"36\u00f72=\", newText: \"12\u00f76=\" },\n  { row: 4, col: 4, oldText: \"70\u00f76=\", newText: \"77\u00f74=\" },\n  // row 8\n  { row: 8, col: 0, oldText: \"77\u00f72=\", newText: \"36\u00f72=\" },\n  { row: 8, col: 1, oldText: \"77\u00f72=\", newText: \"63\u00f77=\" },\n  { row: 8, col: 2, oldText: \"87\u00f76=\", newText: \"73\u00f79=\" },\n  { row: 8, col: 3, oldText: \"78\u00f77=\", newText: \"60\u00f75=\" },\n  { row: 8, col: 4, oldText: \"40\u00f77=\", newText: \"96\u00f74=\" },\n  // row 12\n  { row: 12, col: 0, oldText: \"28\u00f77=\", newText: \"63\u00f79=\" },\n  { row: 12, col: 1, oldText: \"80\u00f74=\", newText: \"83\u00f72=\" },\n  { row: 12, col: 2, oldText: \"60\u00f75=\", newText: \"82\u00f79=\" },\n  { row: 12, col: 3, oldText: \"40\u00f76=\", newText: \"58\u00f74=\" },\n  { row: 12, col: 4, oldText: \"90\u00f78=\", newText: \"77\u00f73=\" },\n  // row 16\n  { row: 16, col: 0, oldText: \"13\u00f77=\", newText: \"87\u00f72=\" },\n  { row: 16, col: 1, oldText: \"29\u00f74=\", newText: \"40\u00f74=\" },\n  { row: 16, col: 2, oldText: \"25\u00f79=\", newText: \"45\u00f79=\" },\n  { row: 16, col: 3, oldText: \"79\u00f75=\", newText: \"29\u00f77=\" },\n  { row: 16, col: 4, oldText: \"67\u00f74=\", newText: \"36\u00f76=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table of division problems in the document.\");\n}\nconst table = tables.items[0];\n\nfor (const { row, col, oldText, newText } of cellReplacements) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Cell (${row}, ${col}) did not contain expected text \"${oldText}\"`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 two-digit division prompts in the practice table with a\n# new set of problems, preserving each cell's existing run/paragraph\n# formatting (font, size, justification).\n#\n# The prompts live in table rows 1, 5, 9, 13 and 17 (Word's 1-based row\n# numbering; 5 cells each, one problem per cell) - the other rows are blank\n# \"work space\" rows. Some old values repeat (\"77\u00f72=\" appears twice) and\n# some new values equal an old value used elsewhere (\"36\u00f72=\", \"60\u00f75=\"), so\n# each replacement is targeted at one specific table cell (by row/column)\n# rather than a document-wide find/replace - that keeps the edit\n# unambiguous and independent of execution order.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$cellReplacements = @(\n    @{ Row = 1;  Col = 1; Old = \"79\u00f72=\"; New = \"22\u00f77=\" },\n    @{ Row = 1;  Col = 2; Old = \"56\u00f78=\"; New = \"62\u00f72=\" },\n    @{ Row = 1;  Col = 3; Old = \"55\u00f78=\"; New = \"28\u00f79=\" },\n    @{ Row = 1;  Col = 4; Old = \"43\u00f72=\"; New = \"73\u00f79=\" },\n    @{ Row = 1;  Col = 5; Old = \"84\u00f74=\"; New = \"40\u00f73=\" },\n\n    @{ Row = 5;  Col = 1; Old = \"33\u00f78=\"; New = \"30\u00f74=\" },\n    @{ Row = 5;  Col = 2; Old = \"97\u00f73=\"; New = \"57\u00f72=\" },\n    @{ Row = 5;  Col = 3; Old = \"46\u00f74=\"; New = \"95\u00f75=\" },\n    @{ Row = 5;  Col = 4; Old = \"36\u00f72=\"; New = \"12\u00f76=\" },\n    @{ Row = 5;  Col = 5; Old = \"70\u00f76=\"; New = \"77\u00f74=\" },\n\n    @{ Row = 9;  Col = 1; Old = \"77\u00f72=\"; New = \"36\u00f72=\" },\n    @{ Row = 9;  Col = 2; Old = \"77\u00f72=\"; New = \"63\u00f77=\" },\n    @{ Row = 9;  Col = 3; Old = \"87\u00f76=\"; New = \"73\u00f79=\" },\n    @{ Row = 9;  Col = 4; Old = \"78\u00f77=\"; New = \"60\u00f75=\" },\n    @{ Row = 9;  Col = 5; Old = \"40\u00f77=\"; New = \"96\u00f74=\" },\n\n    @{ Row = 13; Col = 1; Old = \"28\u00f77=\"; New = \"63\u00f79=\" },\n    @{ Row = 13; Col = 2; Old = \"80\u00f74=\"; New = \"83\u00f72=\" },\n    @{ Row = 13; Col = 3; Old = \"60\u00f75=\"; New = \"82\u00f79=\" },\n    @{ Row = 13; Col = 4; Old = \"40\u00f76=\"; New = \"58\u00f74=\" },\n    @{ Row = 13; Col = 5; Old = \"90\u00f78=\"; New = \"77\u00f73=\" },\n\n    @{ Row = 17; Col = 1; Old = \"13\u00f77=\"; New = \"87\u00f72=\" },\n    @{ Row = 17; Col = 2; Old = \"29\u00f74=\"; New = \"40\u00f74=\" },\n    @{ Row = 17; Col = 3; Old = \"25\u00f79=\"; New = \"45\u00f79=\" },\n    @{ Row = 17; Col = 4; Old = \"79\u00f75=\"; New = \"29\u00f77=\" },\n    @{ Row = 17; Col = 5; Old = \"67\u00f74=\"; New = \"36\u00f76=\" }\n)\n\nforeach ($r in $cellReplacements) {\n    $cell = $t.Cell($r.Row, $r.Col)\n    $range = $cell.Range\n    # Trim the trailing end-of-cell marker from the cell's Range so only\n    # the visible text (\"79\u00f72=\" etc.) is replaced; assigning straight to\n    # Range.Text keeps the run/paragraph formatting (font, size,\n    # justification) untouched and - unlike Find.Execute's Replace, which\n    # this host applies against the first match in the whole story instead\n    # of the given Range - guarantees the edit stays inside this cell only.\n    $textRange = $range.Duplicate\n    $textRange.End = $range.End - 1\n    $textRange.Text = $r.New\n}\n"}
